## Adds a new "2023" year column (T) to the disasters-deaths table, mirroring
## the layout/formatting already used for the "2022" column (S), and widens
## columns A:C to a single uniform width.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths: A, B, C become one uniform width (was 3 distinct widths) ---
$ws.Range("A:C").ColumnWidth = 34.25

# --- New column T: copy the visual formatting of column S (rows 3-34), ---
# --- then overwrite the values with the 2023 data.                      ---
$ws.Range("S3:S34").Copy()
$ws.Range("T3:T34").PasteSpecial(-4122) | Out-Null

$ws.Range("T4").Value  = 2023

$ws.Range("T5").Value  = 44
$ws.Range("T6").Value  = 24
$ws.Range("T7").Value  = 20

$ws.Range("T8").Value  = "-"
$ws.Range("T9").Value  = "-"
$ws.Range("T10").Value = "-"

$ws.Range("T11").Value = 5
$ws.Range("T12").Value = 1
$ws.Range("T13").Value = 4

$ws.Range("T14").Value = 8
$ws.Range("T15").Value = 6
$ws.Range("T16").Value = 2

$ws.Range("T17").Value = 5
$ws.Range("T18").Value = 1
$ws.Range("T19").Value = 4

$ws.Range("T20").Value = 7
$ws.Range("T21").Value = 5
$ws.Range("T22").Value = 2

$ws.Range("T23").Value = "-"
$ws.Range("T24").Value = "-"
$ws.Range("T25").Value = "-"

$ws.Range("T26").Value = 18
$ws.Range("T27").Value = 10
$ws.Range("T28").Value = 8

$ws.Range("T29").Value = "-"
$ws.Range("T30").Value = "-"
$ws.Range("T31").Value = "-"

$ws.Range("T32").Value = 1
$ws.Range("T33").Value = 1
$ws.Range("T34").Value = "-"
